# Update the "Houses of Worship" (column H) indicator values from 1 to 0
# for rows 31 through 176 on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H31:H176").Value = 0
